# ---------------------------------------------------------------------------
# feat: add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" right after "总计" (becomes the 2nd tab),
#    pushing "2022-Q2" / "2021-Q4" / "2021-Q1" / "2020-Q4" one slot later.
# 2. Populate "2022-Q4" with the two new fund rows (mirrors the layout used
#    by the other quarterly sheets).
# 3. Update the "总计" (summary) sheet: insert a new 2022-Q4 row right under
#    the header, shifting the existing rows down by one and renumbering the
#    index column.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$q2 = $wb.Worksheets.Item("2022-Q2")

$newSheet = $wb.Worksheets.Add($null, $total)
$newSheet.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# Step 2: fill in the new sheet with the same layout as the other quarterly
# fund-holding sheets (header row + fund rows).
# ---------------------------------------------------------------------------
function Set-HeaderCell($ws, $addr, $text) {
    $ws.Range($addr).Value = $text
    $ws.Range($addr).Font.Bold = $true
    $ws.Range($addr).Borders.LineStyle = 1
    $ws.Range($addr).HorizontalAlignment = -4108
    $ws.Range($addr).VerticalAlignment = -4160
}

function Set-IndexCell($ws, $addr, $num) {
    $ws.Range($addr).Value = $num
    $ws.Range($addr).Font.Bold = $true
    $ws.Range($addr).Borders.LineStyle = 1
    $ws.Range($addr).HorizontalAlignment = -4108
    $ws.Range($addr).VerticalAlignment = -4160
}

function Set-TextCell($ws, $addr, $text) {
    $ws.Range($addr).Value = "'" + $text
}

function Set-NumCell($ws, $addr, $num) {
    $ws.Range($addr).Value = $num
}

Set-HeaderCell $newSheet "B1" "基金代码"
Set-HeaderCell $newSheet "C1" "基金名称"
Set-HeaderCell $newSheet "D1" "基金规模"
Set-HeaderCell $newSheet "E1" "股票总仓位"
Set-HeaderCell $newSheet "F1" "仓位占比"
Set-HeaderCell $newSheet "G1" "持有市值(亿元)"
Set-HeaderCell $newSheet "H1" "仓位排名"

Set-IndexCell $newSheet "A2" 0
Set-TextCell  $newSheet "B2" "010377"
Set-TextCell  $newSheet "C2" "广发价值核心混合A"
Set-TextCell  $newSheet "D2" "24.90"
Set-TextCell  $newSheet "E2" "89.01"
Set-TextCell  $newSheet "F2" "5.58"
Set-TextCell  $newSheet "G2" "1.3894"
Set-NumCell   $newSheet "H2" 3

Set-IndexCell $newSheet "A3" 1
Set-TextCell  $newSheet "B3" "010378"
Set-TextCell  $newSheet "C3" "广发价值核心混合C"
Set-TextCell  $newSheet "D3" "4.47"
Set-TextCell  $newSheet "E3" "89.01"
Set-TextCell  $newSheet "F3" "5.58"
Set-TextCell  $newSheet "G3" "0.2494"
Set-NumCell   $newSheet "H3" 3

# ---------------------------------------------------------------------------
# Step 3: update the "总计" sheet — insert the 2022-Q4 summary row and shift
# the remaining rows down.
# ---------------------------------------------------------------------------
$totals = @(
    @{ Label = "2022-Q4"; Count = 2; Value = 1.64 },
    @{ Label = "2022-Q2"; Count = 2; Value = 0.04 },
    @{ Label = "2021-Q4"; Count = 2; Value = 0.25 },
    @{ Label = "2021-Q1"; Count = 6; Value = 4.79 },
    @{ Label = "2020-Q4"; Count = 3; Value = 2.31 }
)

for ($i = 0; $i -lt $totals.Count; $i++) {
    $row = $i + 2
    $entry = $totals[$i]
    $total.Cells.Item($row, 1).Value = $i
    $total.Cells.Item($row, 1).Font.Bold = $true
    $total.Cells.Item($row, 1).Borders.LineStyle = 1
    $total.Cells.Item($row, 1).HorizontalAlignment = -4108
    $total.Cells.Item($row, 1).VerticalAlignment = -4160
    $total.Cells.Item($row, 2).Value = $entry.Label
    $total.Cells.Item($row, 3).Value = $entry.Count
    $total.Cells.Item($row, 4).Value = $entry.Value
}

# ---------------------------------------------------------------------------
# Restore the originally-active tab (last sheet, "2020-Q4").
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
